$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D2:E51 so numeric-looking strings (e.g. "1.001") are
# not auto-converted to numbers by Excel's smart-entry parsing, matching the
# original inline-string cell typing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.058.16"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.762.68"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "335.47"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.3903"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("D8").Value = "0.3406"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "45.39"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").Value = "1.125"
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("D11").Value = "0.07237"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "22.35"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").Value = "6.159"
$ws.Range("E14").Value = "  -4.56%  "
$ws.Range("D15").Value = "1.757.30"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "7.088"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").Value = "0.00001059"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "0.06626"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "80.61"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "0.9981"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "16.96"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "6.217"
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("D23").Value = "28.030.01"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").Value = "11.66"
$ws.Range("E24").Value = "  -3.29%  "
$ws.Range("D25").Value = "2.393"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "154.67"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "19.98"
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("D28").Value = "2.307"
$ws.Range("E28").Value = "  -4.25%  "
$ws.Range("D29").Value = "1.957.70"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "1.283"
$ws.Range("E30").Value = "  -11.85%  "
$ws.Range("D31").Value = "129.24"
$ws.Range("E31").Value = "  -5.08%  "
$ws.Range("D32").Value = "4.078"
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("D33").Value = "5.832"
$ws.Range("E33").Value = "  -4.60%  "
$ws.Range("D34").Value = "0.08699"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").Value = "12.10"
$ws.Range("E35").Value = "  -5.20%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "5.146"
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.06165"
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02285"
$ws.Range("E38").Value = "  -6.15%  "
$ws.Range("D39").Value = "0.6486"
$ws.Range("E39").Value = "  -5.27%  "
$ws.Range("D40").Value = "0.2112"
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("D41").Value = "1.501"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D43").Value = "0.9983"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "7.891"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("D45").Value = "13.75"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").Value = "3.829"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "0.5997"
$ws.Range("E47").Value = "  -4.71%  "
$ws.Range("D48").Value = "126.71"
$ws.Range("E48").Value = "  -4.92%  "
$ws.Range("D49").Value = "1.980"
$ws.Range("E49").Value = "  -5.01%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.157"
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.07001"
$ws.Range("E51").Value = "  -5.89%  "

# Restore default (unstyled) appearance now that the text is committed.
$ws.Range("D2:E51").Style = "Normal"
